# Mise a jour des documents administratifs
# Fill in the journal entries for 17.05.2022 (rows 58-61) that were left
# blank, and move the sheet's active selection forward to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 58 : 10:27 - 12:15 -> documentation
$ws.Range("A58").Value = 44698
$ws.Range("B58").Value = 0.43541666666666662
$ws.Range("C58").Value = 0.51041666666666663
$ws.Range("E58").Value = "documentation"

# Row 59 : 13:30 - 14:50 -> HPA: Conception
$ws.Range("A59").Value = 44698
$ws.Range("B59").Value = 0.5625
$ws.Range("C59").Value = 0.61805555555555558
$ws.Range("E59").Value = "HPA: Conception"

# Row 60 : 14:50 - 15:30 -> HPA: Conception
$ws.Range("A60").Value = 44698
$ws.Range("B60").Value = 0.61805555555555558
$ws.Range("C60").Value = 0.64583333333333337
$ws.Range("E60").Value = "HPA: Conception"

# Row 61 : 15:30 - 16:05 -> Entretien avec le second expert
$ws.Range("A61").Value = 44698
$ws.Range("B61").Value = 0.64583333333333337
$ws.Range("C61").Value = 0.67013888888888884
$ws.Range("E61").Value = "Entretien avec le second expert"

# Scroll the view down a bit and move the selection, matching the author's
# on-screen state after typing the new rows.
$win = $excel.ActiveWindow
$win.ScrollRow = 56
$win.ScrollColumn = 1
$null = $ws.Range("E62").Select()
